$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 439
$col = 3  # Column C

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, $col).Value = 45177
}
